$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing last row (row 140) values per diff
$ws.Cells.Item(140, 2).Value = 6.6   # B140: 6.5 -> 6.6
$ws.Cells.Item(140, 4).Value = 6.7   # D140: 6.6 -> 6.7

# Append a new data row (row 141) with a new date label and values.
# Use a scratch cell to enter the date-looking text as plain text
# (NumberFormat "@") so Excel doesn't auto-convert it to a date serial,
# then copy just the value into place and clean up the scratch cell.
$scratch = $ws.Cells.Item(500, 1)
$scratch.NumberFormat = "@"
$scratch.Value = "01-08-2021"
$scratch.Copy()
$ws.Cells.Item(141, 1).PasteSpecial(-4163)
$scratch.Clear()

$ws.Cells.Item(141, 2).Value = 6.4
$ws.Cells.Item(141, 3).Value = 1.6
$ws.Cells.Item(141, 4).Value = 6.7
